# Add a new "forebrain_neurons" column to the human data sheet, populated
# with the (human) forebrain neuron count for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work out where the data currently ends and which column is free.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$newCol = $usedRange.Column + $usedRange.Columns.Count

# New header for the new column.
$ws.Cells.Item(1, $newCol).Value = "forebrain_neurons"

# Every data row (everything below the header) gets the same constant:
# humans' forebrain neuron count.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol).Value = 24560000000
}
